$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.189.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.446.97"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.69%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.57"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "187.42"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.72%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.437.83"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.73%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.640"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "57.60"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +7.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000275"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.40"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.997.40"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.64%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.447.95"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.076.37"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.58%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.03"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.85%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "487.12"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.61"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +14.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.19"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +20.62%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.31"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.34"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.92%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.91"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.00%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.18"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.32"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +11.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "605.01"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "64.79"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.78"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.111"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.90%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.145"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.88"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0776"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.385"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.14%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.186.95"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.87"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0428"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.82%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.56"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.23"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.65"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +15.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.06%  "

# Row 50 becomes Monero, Row 51 becomes THORChain (swap with updated values)
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.64"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.14%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.59"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.23%  "
